$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataCombined")

$oldValue = "Laskin 1982.Group A_Aciclovir_1_Human__PeripheralVenousBlood_Plasma_2.5 mg/kg_iv_"
$newValue = "Laskin 1982.Group A_Aciclovir_1_Human_MALE_PeripheralVenousBlood_Plasma_2.5 mg/kg_iv_"

$usedRange = $ws.UsedRange
foreach ($cell in $usedRange.Cells) {
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
